# Insert a new row at position 219. Excel shifts rows 219:245 down to 220:246,
# carrying their formatting (including the date style on column D) with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new record's data.
$ws.Range("A219").Value = 11
$ws.Range("B219").Value = "Vega Monumental Concepción"
$ws.Range("C219").Value = "Bíobío"
$ws.Range("D219").Value2 = 44946
$ws.Range("E219").Value = 8
$ws.Range("F219").Value = 100112003
$ws.Range("G219").Value = "Ajo"
$ws.Range("H219").Value = "Chino"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 220
$ws.Range("K219").Value = 14000
$ws.Range("L219").Value = 15000
$ws.Range("M219").Value = 14455
$ws.Range("N219").Value = "$/caja 10 kilos"
$ws.Range("O219").Value = "China"
$ws.Range("P219").Value = 1446
$ws.Range("Q219").Value = 10
$ws.Range("R219").Value = "Hortaliza"
